# Revert "RESTORE: Recover all 973 original multi-industry template files ..."
# i.e. re-apply the Product/AI-ML -> Product wording changes across the
# Staffing Plan workbook.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet: Resource Overview
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Overview")

$ws.Range("A2").Value = "Product Development Implementation Project"
$ws.Range("B6").Value = "Enterprise Product Development Implementation"
$ws.Range("A18").Value = "Data Science/Product"
$ws.Range("G18").Value = "Product, Python, Statistics"

# incidental empty-row touch (row 13 was visited while the original
# template was regenerated, leaving a blank row behind)
$ws.Rows.Item(13).RowHeight = 15

# ------------------------------------------------------------------
# Sheet: Detailed Staffing Plan
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Detailed Staffing Plan")

$ws.Range("A1").Value = "DETProductLED STAFFING PLAN"
$ws.Rows.Item(2).RowHeight = 15

$ws.Range("C9").Value = "Data Science/Product"
$ws.Range("K9").Value = "Product, Deep Learning, Python"
$ws.Range("P9").Value = "Product Lead"

$ws.Range("C10").Value = "Data Science/Product"
$ws.Range("K10").Value = "Product, Statistics, R/Python"

$ws.Range("C11").Value = "Data Science/Product"
$ws.Range("K11").Value = "Product, Python, Visualization"

$ws.Range("B12").Value = "Product Engineer"
$ws.Range("C12").Value = "Data Science/Product"
$ws.Range("K12").Value = "ProductOps, Python, Cloud"

$ws.Range("C13").Value = "Data Science/Product"

# ------------------------------------------------------------------
# Sheet: Resource Timeline
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Timeline")

$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15

# ------------------------------------------------------------------
# Sheet: Skills Matrix
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Skills Matrix")

$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15
$ws.Range("D3").Value = "Product Innovation"

# ------------------------------------------------------------------
# Sheet: Cost Analysis
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cost Analysis")

$ws.Rows.Item(2).RowHeight = 15
$ws.Range("A6").Value = "Data Science/Product"
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 15

# ------------------------------------------------------------------
# Sheet: Resource Risk Assessment
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resource Risk Assessment")

$ws.Rows.Item(2).RowHeight = 15
$ws.Range("B5").Value = "Team lacks required Product expertise"
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
